# Swap the contents of columns B:G between the given pairs of rows.
# These pairs represent two stock-report lines that were recorded in the
# wrong order (the second scan's data ended up on the first row and vice
# versa); this restores the correct row <-> data association.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(142, 143),
    @(256, 257),
    @(305, 306),
    @(338, 339),
    @(342, 344),
    @(364, 365),
    @(374, 375),
    @(381, 382),
    @(392, 393),
    @(411, 412),
    @(413, 414),
    @(423, 424),
    @(449, 450),
    @(578, 579),
    @(596, 597),
    @(679, 680),
    @(701, 702),
    @(712, 713),
    @(864, 865)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:G$r1")
    $range2 = $ws.Range("B$r2`:G$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value = $values2
    $range2.Value = $values1
}
